# Applies data-correction edits to the team game-log worksheet:
# - Updates ~657 statistic cells that were off by one game (GP/W/etc. and
#   their dependent rates/ranks) due to stats being pulled a day early.
# - Normalizes the Date column (BF) from "5-2-2009-10" to ISO "2010-05-02",
#   forcing the cells to remain plain text (not auto-converted to a date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cell updates: @(row, column, newValue)
$numericUpdates = @(
    @(2, 4, 82),
    @(2, 5, 53),
    @(2, 7, 0.646),
    @(2, 9, 38.8),
    @(2, 11, 0.468),
    @(2, 13, 17.7),
    @(2, 14, 0.36),
    @(2, 18, 11.8),
    @(2, 20, 41.7),
    @(2, 21, 21.8),
    @(2, 23, 7.2),
    @(2, 26, 19.9),
    @(2, 27, 19.3),
    @(2, 28, 101.7),
    @(2, 29, 4.7),
    @(2, 30, 1),
    @(2, 31, 6),
    @(2, 32, 6),
    @(2, 33, 6),
    @(2, 35, 6),
    @(2, 40, 9),
    @(2, 41, 21),
    @(2, 44, 6),
    @(2, 45, 23),
    @(2, 46, 17),
    @(2, 49, 15),
    @(2, 52, 8),
    @(2, 54, 13),
    @(3, 4, 82),
    @(3, 5, 50),
    @(3, 7, 0.61),
    @(3, 8, 48.2),
    @(3, 9, 37.1),
    @(3, 11, 0.483),
    @(3, 12, 6.1),
    @(3, 13, 17.5),
    @(3, 14, 0.348),
    @(3, 16, 25.5),
    @(3, 17, 0.746),
    @(3, 18, 8.699999999999999),
    @(3, 19, 29.9),
    @(3, 20, 38.6),
    @(3, 24, 4.9),
    @(3, 28, 99.2),
    @(3, 29, 3.7),
    @(3, 30, 1),
    @(3, 31, 9),
    @(3, 33, 9),
    @(3, 38, 16),
    @(3, 39, 16),
    @(3, 40, 17),
    @(3, 41, 15),
    @(3, 43, 21),
    @(3, 50, 15),
    @(3, 51, 15),
    @(3, 52, 22),
    @(3, 55, 9),
    @(4, 4, 82),
    @(4, 6, 38),
    @(4, 7, 0.537),
    @(4, 9, 34.9),
    @(4, 11, 0.453),
    @(4, 12, 5.6),
    @(4, 13, 16.2),
    @(4, 14, 0.346),
    @(4, 17, 0.751),
    @(4, 18, 10.5),
    @(4, 19, 30.4),
    @(4, 21, 20.2),
    @(4, 24, 5.4),
    @(4, 26, 19.5),
    @(4, 28, 95.3),
    @(4, 29, 1.5),
    @(4, 30, 1),
    @(4, 37, 22),
    @(4, 40, 21),
    @(4, 42, 5),
    @(4, 43, 20),
    @(4, 44, 21),
    @(4, 50, 7),
    @(4, 51, 29),
    @(4, 52, 5),
    @(4, 54, 28),
    @(4, 55, 15),
    @(5, 31, 17),
    @(5, 34, 2),
    @(5, 36, 10),
    @(5, 48, 15),
    @(5, 49, 24),
    @(5, 51, 23),
    @(6, 4, 82),
    @(6, 6, 21),
    @(6, 7, 0.744),
    @(6, 8, 48.2),
    @(6, 9, 37.8),
    @(6, 10, 77.90000000000001),
    @(6, 11, 0.485),
    @(6, 14, 0.381),
    @(6, 15, 19.1),
    @(6, 16, 26.6),
    @(6, 18, 9.6),
    @(6, 19, 32.8),
    @(6, 21, 22.4),
    @(6, 25, 4),
    @(6, 26, 19.4),
    @(6, 27, 20.8),
    @(6, 28, 102.1),
    @(6, 29, 6.5),
    @(6, 30, 1),
    @(6, 35, 15),
    @(6, 36, 28),
    @(6, 41, 12),
    @(6, 42, 5),
    @(6, 47, 6),
    @(6, 48, 12),
    @(6, 53, 15),
    @(7, 4, 82),
    @(7, 6, 27),
    @(7, 7, 0.671),
    @(7, 9, 38.3),
    @(7, 10, 82.40000000000001),
    @(7, 11, 0.464),
    @(7, 14, 0.372),
    @(7, 15, 18.6),
    @(7, 16, 22.8),
    @(7, 21, 23.4),
    @(7, 22, 12.9),
    @(7, 24, 5.5),
    @(7, 26, 19.1),
    @(7, 27, 20.1),
    @(7, 28, 102),
    @(7, 29, 2.7),
    @(7, 30, 1),
    @(7, 35, 11),
    @(7, 40, 5),
    @(7, 42, 25),
    @(7, 50, 6),
    @(7, 53, 21),
    @(8, 4, 82),
    @(8, 5, 53),
    @(8, 9, 38.1),
    @(8, 10, 81.40000000000001),
    @(8, 12, 6.6),
    @(8, 13, 18.5),
    @(8, 14, 0.359),
    @(8, 19, 30.5),
    @(8, 23, 8.300000000000001),
    @(8, 25, 5.3),
    @(8, 28, 106.5),
    @(8, 29, 4.1),
    @(8, 30, 1),
    @(8, 31, 6),
    @(8, 32, 6),
    @(8, 33, 6),
    @(8, 34, 17),
    @(8, 35, 13),
    @(8, 40, 10),
    @(8, 44, 17),
    @(8, 45, 16),
    @(8, 48, 11),
    @(8, 51, 23),
    @(9, 4, 82),
    @(9, 5, 27),
    @(9, 7, 0.329),
    @(9, 9, 35.9),
    @(9, 10, 80.5),
    @(9, 12, 4.6),
    @(9, 14, 0.314),
    @(9, 15, 17.7),
    @(9, 16, 24.4),
    @(9, 17, 0.728),
    @(9, 18, 12.8),
    @(9, 22, 13.4),
    @(9, 27, 20.8),
    @(9, 29, -5.1),
    @(9, 30, 1),
    @(9, 31, 24),
    @(9, 32, 24),
    @(9, 33, 24),
    @(9, 34, 17),
    @(9, 39, 26),
    @(9, 46, 27),
    @(9, 48, 7),
    @(9, 49, 14),
    @(9, 53, 16),
    @(10, 4, 82),
    @(10, 6, 56),
    @(10, 7, 0.317),
    @(10, 12, 7.7),
    @(10, 13, 20.6),
    @(10, 15, 19.9),
    @(10, 16, 25.4),
    @(10, 17, 0.782),
    @(10, 18, 9.199999999999999),
    @(10, 25, 5),
    @(10, 30, 1),
    @(10, 31, 26),
    @(10, 32, 26),
    @(10, 33, 26),
    @(10, 40, 4),
    @(10, 43, 5),
    @(10, 47, 5),
    @(10, 55, 22),
    @(11, 4, 82),
    @(11, 5, 42),
    @(11, 6, 40),
    @(11, 7, 0.512),
    @(11, 8, 48.5),
    @(11, 9, 37.7),
    @(11, 10, 84.40000000000001),
    @(11, 11, 0.447),
    @(11, 18, 11.8),
    @(11, 19, 30.1),
    @(11, 20, 42),
    @(11, 22, 14.5),
    @(11, 23, 7.1),
    @(11, 25, 6.5),
    @(11, 26, 20.9),
    @(11, 28, 102.4),
    @(11, 29, -0.4),
    @(11, 30, 1),
    @(11, 38, 5),
    @(11, 41, 14),
    @(11, 45, 21),
    @(11, 49, 18),
    @(11, 52, 17),
    @(12, 4, 82),
    @(12, 6, 50),
    @(12, 7, 0.39),
    @(12, 10, 83.2),
    @(12, 11, 0.443),
    @(12, 13, 23.1),
    @(12, 14, 0.348),
    @(12, 15, 19.1),
    @(12, 17, 0.775),
    @(12, 22, 15),
    @(12, 23, 7.1),
    @(12, 28, 100.8),
    @(12, 29, -3),
    @(12, 30, 1),
    @(12, 34, 29),
    @(12, 36, 11),
    @(12, 40, 18),
    @(12, 41, 13),
    @(12, 45, 6),
    @(12, 47, 15),
    @(12, 48, 25),
    @(12, 49, 17),
    @(12, 51, 18),
    @(13, 4, 82),
    @(13, 6, 53),
    @(13, 7, 0.354),
    @(13, 9, 36.6),
    @(13, 11, 0.455),
    @(13, 13, 17.8),
    @(13, 15, 16.6),
    @(13, 17, 0.73),
    @(13, 21, 22.1),
    @(13, 22, 15.7),
    @(13, 23, 6.5),
    @(13, 26, 19.3),
    @(13, 28, 95.7),
    @(13, 30, 1),
    @(13, 37, 20),
    @(13, 41, 28),
    @(13, 42, 26),
    @(13, 45, 17),
    @(13, 46, 14),
    @(13, 47, 9),
    @(13, 49, 23),
    @(13, 51, 9),
    @(13, 54, 27),
    @(14, 4, 82),
    @(14, 5, 57),
    @(14, 7, 0.695),
    @(14, 9, 38.3),
    @(14, 11, 0.457),
    @(14, 13, 19),
    @(14, 14, 0.341),
    @(14, 15, 18.5),
    @(14, 16, 24.2),
    @(14, 17, 0.765),
    @(14, 18, 11.9),
    @(14, 22, 13.4),
    @(14, 23, 7.5),
    @(14, 26, 19.4),
    @(14, 27, 21.2),
    @(14, 30, 1),
    @(14, 35, 9),
    @(14, 37, 18),
    @(14, 42, 17),
    @(14, 44, 4),
    @(14, 50, 16),
    @(14, 51, 9),
    @(14, 52, 4),
    @(15, 4, 82),
    @(15, 6, 42),
    @(15, 7, 0.488),
    @(15, 8, 48.5),
    @(15, 9, 39.3),
    @(15, 10, 83.8),
    @(15, 11, 0.469),
    @(15, 14, 0.337),
    @(15, 16, 26.9),
    @(15, 17, 0.733),
    @(15, 20, 43.5),
    @(15, 21, 18.8),
    @(15, 26, 20.2),
    @(15, 28, 102.5),
    @(15, 29, -1.5),
    @(15, 30, 1),
    @(15, 32, 18),
    @(15, 34, 2),
    @(15, 35, 4),
    @(15, 36, 7),
    @(15, 37, 9),
    @(15, 40, 26),
    @(15, 41, 8),
    @(15, 47, 29),
    @(15, 50, 18),
    @(15, 51, 28),
    @(15, 53, 2),
    @(15, 54, 7),
    @(16, 4, 82),
    @(16, 5, 47),
    @(16, 7, 0.573),
    @(16, 9, 36.4),
    @(16, 10, 79.5),
    @(16, 11, 0.458),
    @(16, 14, 0.346),
    @(16, 15, 17.7),
    @(16, 16, 23.5),
    @(16, 17, 0.752),
    @(16, 22, 13.2),
    @(16, 28, 96.5),
    @(16, 29, 2.3),
    @(16, 30, 1),
    @(16, 35, 26),
    @(16, 37, 17),
    @(16, 38, 17),
    @(16, 39, 17),
    @(16, 40, 19),
    @(16, 43, 19),
    @(16, 46, 15),
    @(16, 47, 28),
    @(16, 49, 12),
    @(16, 52, 16),
    @(17, 31, 14),
    @(17, 38, 6),
    @(17, 49, 19),
    @(17, 50, 20),
    @(17, 54, 23),
    @(17, 55, 14),
    @(18, 4, 82),
    @(18, 5, 15),
    @(18, 7, 0.183),
    @(18, 10, 84.40000000000001),
    @(18, 11, 0.449),
    @(18, 12, 4.9),
    @(18, 13, 14.4),
    @(18, 14, 0.341),
    @(18, 15, 17.5),
    @(18, 16, 23.5),
    @(18, 20, 42.9),
    @(18, 21, 19.8),
    @(18, 23, 7.3),
    @(18, 24, 3.7),
    @(18, 25, 5.4),
    @(18, 27, 20.6),
    @(18, 29, -9.6),
    @(18, 30, 1),
    @(18, 34, 17),
    @(18, 35, 14),
    @(18, 36, 3),
    @(18, 37, 25),
    @(18, 39, 28),
    @(18, 43, 22),
    @(18, 44, 9),
    @(18, 49, 13),
    @(18, 51, 25),
    @(18, 53, 17),
    @(19, 4, 82),
    @(19, 6, 70),
    @(19, 7, 0.146),
    @(19, 12, 4.6),
    @(19, 14, 0.318),
    @(19, 15, 19.2),
    @(19, 16, 24.6),
    @(19, 17, 0.78),
    @(19, 19, 28.8),
    @(19, 20, 39.7),
    @(19, 22, 14.4),
    @(19, 24, 4.8),
    @(19, 25, 5.1),
    @(19, 27, 20.1),
    @(19, 29, -9.1),
    @(19, 30, 1),
    @(19, 34, 17),
    @(19, 39, 27),
    @(19, 43, 6),
    @(19, 48, 16),
    @(19, 49, 20),
    @(19, 51, 21),
    @(19, 52, 10),
    @(19, 53, 23),
    @(20, 4, 82),
    @(20, 6, 45),
    @(20, 7, 0.451),
    @(20, 10, 83.40000000000001),
    @(20, 13, 19.2),
    @(20, 14, 0.363),
    @(20, 16, 20.3),
    @(20, 17, 0.778),
    @(20, 18, 10.4),
    @(20, 20, 40.3),
    @(20, 21, 22.3),
    @(20, 23, 7.6),
    @(20, 27, 19.5),
    @(20, 28, 100.2),
    @(20, 29, -2.5),
    @(20, 30, 1),
    @(20, 35, 7),
    @(20, 40, 8),
    @(20, 43, 7),
    @(20, 44, 22),
    @(20, 45, 24),
    @(20, 46, 25),
    @(20, 47, 8),
    @(20, 51, 12),
    @(21, 4, 82),
    @(21, 6, 53),
    @(21, 7, 0.354),
    @(21, 9, 38.1),
    @(21, 10, 83.90000000000001),
    @(21, 11, 0.455),
    @(21, 12, 9.1),
    @(21, 13, 26.2),
    @(21, 14, 0.346),
    @(21, 28, 102.1),
    @(21, 29, -3.8),
    @(21, 30, 1),
    @(21, 32, 22),
    @(21, 33, 22),
    @(21, 36, 6),
    @(21, 37, 21),
    @(21, 40, 20),
    @(21, 43, 4),
    @(21, 45, 20),
    @(21, 51, 13),
    @(21, 52, 9),
    @(21, 53, 29),
    @(21, 54, 9),
    @(21, 55, 23),
    @(22, 4, 82),
    @(22, 5, 50),
    @(22, 7, 0.61),
    @(22, 9, 37.4),
    @(22, 10, 80.8),
    @(22, 11, 0.462),
    @(22, 14, 0.34),
    @(22, 15, 21.7),
    @(22, 18, 11.7),
    @(22, 29, 3.5),
    @(22, 30, 1),
    @(22, 40, 25),
    @(22, 44, 10),
    @(22, 51, 14),
    @(22, 52, 20),
    @(22, 54, 14),
    @(22, 55, 10),
    @(23, 4, 82),
    @(23, 5, 59),
    @(23, 7, 0.72),
    @(23, 9, 36.6),
    @(23, 10, 78),
    @(23, 11, 0.47),
    @(23, 12, 10.3),
    @(23, 14, 0.375),
    @(23, 15, 19.2),
    @(23, 16, 26.5),
    @(23, 17, 0.724),
    @(23, 20, 43.2),
    @(23, 21, 19.7),
    @(23, 22, 14.1),
    @(23, 23, 6.2),
    @(23, 28, 102.8),
    @(23, 29, 7.5),
    @(23, 30, 1),
    @(23, 36, 27),
    @(23, 37, 7),
    @(23, 40, 3),
    @(23, 42, 7),
    @(23, 52, 7),
    @(23, 54, 6),
    @(24, 4, 82),
    @(24, 6, 55),
    @(24, 7, 0.329),
    @(24, 15, 16.7),
    @(24, 19, 29.5),
    @(24, 20, 41),
    @(24, 22, 14.5),
    @(24, 24, 5.4),
    @(24, 26, 20.5),
    @(24, 28, 97.7),
    @(24, 29, -3.9),
    @(24, 30, 1),
    @(24, 41, 27),
    @(24, 50, 8),
    @(24, 54, 22),
    @(25, 4, 82),
    @(25, 5, 54),
    @(25, 7, 0.659),
    @(25, 13, 21.6),
    @(25, 17, 0.77),
    @(25, 19, 31.9),
    @(25, 20, 43),
    @(25, 21, 23.3),
    @(25, 25, 4.5),
    @(25, 30, 1),
    @(25, 41, 5),
    @(25, 42, 9),
    @(25, 45, 7),
    @(25, 51, 11),
    @(25, 52, 18),
    @(25, 55, 5),
    @(26, 4, 82),
    @(26, 5, 50),
    @(26, 7, 0.61),
    @(26, 12, 6),
    @(26, 14, 0.354),
    @(26, 17, 0.79),
    @(26, 19, 29.1),
    @(26, 20, 40.2),
    @(26, 22, 12.3),
    @(26, 23, 6.4),
    @(26, 24, 4.3),
    @(26, 26, 20.9),
    @(26, 29, 3.3),
    @(26, 30, 1),
    @(26, 35, 27),
    @(26, 40, 13),
    @(26, 46, 26),
    @(26, 52, 19),
    @(27, 4, 82),
    @(27, 6, 57),
    @(27, 7, 0.305),
    @(27, 9, 38.3),
    @(27, 11, 0.456),
    @(27, 14, 0.349),
    @(27, 17, 0.726),
    @(27, 19, 30.7),
    @(27, 20, 42.6),
    @(27, 21, 20.5),
    @(27, 28, 100),
    @(27, 29, -4.4),
    @(27, 30, 1),
    @(27, 31, 28),
    @(27, 32, 28),
    @(27, 33, 28),
    @(27, 35, 9),
    @(27, 37, 19),
    @(27, 39, 20),
    @(27, 42, 18),
    @(27, 48, 23),
    @(27, 49, 21),
    @(27, 53, 22),
    @(28, 4, 82),
    @(28, 5, 50),
    @(28, 7, 0.61),
    @(28, 9, 38.4),
    @(28, 10, 81.2),
    @(28, 11, 0.473),
    @(28, 12, 6.8),
    @(28, 13, 18.9),
    @(28, 16, 24),
    @(28, 17, 0.74),
    @(28, 21, 22.3),
    @(28, 22, 13.6),
    @(28, 24, 4.6),
    @(28, 26, 20.4),
    @(28, 28, 101.4),
    @(28, 29, 5.1),
    @(28, 30, 1),
    @(28, 34, 17),
    @(28, 35, 8),
    @(28, 42, 18),
    @(28, 43, 24),
    @(28, 44, 18),
    @(28, 47, 7),
    @(28, 50, 22),
    @(29, 4, 82),
    @(29, 5, 40),
    @(29, 7, 0.488),
    @(29, 9, 39),
    @(29, 11, 0.482),
    @(29, 15, 19.7),
    @(29, 16, 25.8),
    @(29, 18, 9.800000000000001),
    @(29, 19, 30.6),
    @(29, 20, 40.4),
    @(29, 26, 22.2),
    @(29, 27, 21.1),
    @(29, 29, -1.8),
    @(29, 30, 1),
    @(29, 31, 18),
    @(29, 32, 18),
    @(29, 33, 18),
    @(29, 34, 17),
    @(29, 41, 7),
    @(29, 42, 8),
    @(29, 45, 15),
    @(29, 46, 23),
    @(29, 47, 10),
    @(29, 50, 21),
    @(29, 52, 23),
    @(30, 4, 82),
    @(30, 6, 29),
    @(30, 7, 0.646),
    @(30, 12, 5.4),
    @(30, 13, 14.7),
    @(30, 14, 0.364),
    @(30, 15, 20.2),
    @(30, 17, 0.741),
    @(30, 23, 8.199999999999999),
    @(30, 27, 22.2),
    @(30, 29, 5.3),
    @(30, 30, 1),
    @(30, 31, 6),
    @(30, 32, 6),
    @(30, 33, 6),
    @(30, 35, 3),
    @(30, 40, 7),
    @(30, 43, 23),
    @(30, 50, 16),
    @(31, 4, 82),
    @(31, 5, 26),
    @(31, 7, 0.317),
    @(31, 13, 14.9),
    @(31, 17, 0.762),
    @(31, 19, 30),
    @(31, 20, 41.8),
    @(31, 25, 5.1),
    @(31, 26, 21.4),
    @(31, 29, -4.8),
    @(31, 30, 1),
    @(31, 31, 26),
    @(31, 32, 26),
    @(31, 33, 26),
    @(31, 37, 24),
    @(31, 40, 14),
    @(31, 44, 7),
    @(31, 45, 22),
    @(31, 46, 16),
    @(31, 48, 21),
    @(31, 51, 19),
    @(31, 52, 21),
    @(31, 53, 18)
)

foreach ($u in $numericUpdates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# Force the Date column to remain text (prevent auto date-serial conversion)
# before writing the corrected ISO-formatted date strings.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Every data row (2-31) gets the same corrected Date value.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).Value = "2010-05-02"
}